{"js": "// The author's edit removes the stray duplicated word \" souhv\u011bzd\u00ed \" so\n// that \"...zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed Orion.\" becomes\n// \"...zobrazuj\u00edSouhv\u011bzd\u00ed Orion.\" (run text collapses together, no space).\n// This exact phrase occurs 4 times across the document body (identical\n// paragraphs repeated throughout the activity guide); replace every one.\n\nconst searchText = \"zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\";\nconst replacementText = \"zobrazuj\u00edSouhv\u011bzd\u00ed\";\n\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The author's edit removes the stray duplicated word \" souhv\u011bzd\u00ed \" so that\n# \"...zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed Orion.\" becomes\n# \"...zobrazuj\u00edSouhv\u011bzd\u00ed Orion.\" (the run text collapses together, no space\n# left behind). This exact phrase is repeated 4 times throughout the\n# document body (identical paragraphs recurring through the activity\n# guide) - replace every occurrence.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\"\n$find.Replacement.Text = \"zobrazuj\u00edSouhv\u011bzd\u00ed\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n    [ref]$find.Text,\n    [ref]$find.MatchCase,\n    [ref]$find.MatchWholeWord,\n    [ref]$find.MatchWildcards,\n    $null, $null, [ref]$find.Forward, [ref]$find.Wrap, $null,\n    [ref]$find.Replacement.Text,\n    2  # wdReplaceAll\n) | Out-Null\n"}
